# Update the cryptos list — price (D) and 1h volume (E) columns, plus
# the row 49/50 swap (HuobiToken <-> MultiversX moved position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.916.69"
$ws.Range("E2").Value = "  -0.54%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.302.87"
$ws.Range("E3").Value = "  -0.23%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "305.57"
$ws.Range("E5").Value = "  +1.54%  "

# Row 6 - Solana
$ws.Range("D6").Value = "97.54"
$ws.Range("E6").Value = "  -0.57%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.85%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -2.18%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "35.79"
$ws.Range("E10").Value = "  -0.17%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +0.18%  "

# Row 12 - Chainlink
$ws.Range("E12").Value = "  +0.77%  "

# Row 13 - TRON
$ws.Range("D13").Value = "0.118"
$ws.Range("E13").Value = "  +0.98%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "6.79"
$ws.Range("E14").Value = "  -1.46%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.660.19"
$ws.Range("E15").Value = "  -0.06%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.302.54"
$ws.Range("E16").Value = "  +1.00%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  -0.27%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "42.848.08"
$ws.Range("E18").Value = "  -0.37%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("D19").Value = "12.73"
$ws.Range("E19").Value = "  -4.41%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -0.51%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -1.12%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "67.79"
$ws.Range("E22").Value = "  -1.04%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "236.84"
$ws.Range("E23").Value = "  -0.53%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  -1.45%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +1.68%  "

# Row 26 - Dai (unchanged)

# Row 27 - LEO
$ws.Range("E27").Value = "  -0.13%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "25.46"
$ws.Range("E28").Value = "  +2.50%  "

# Row 29 - Monero
$ws.Range("D29").Value = "167.11"
$ws.Range("E29").Value = "  -0.58%  "

# Row 30 - Toncoin
$ws.Range("D30").Value = "2.06"
$ws.Range("E30").Value = "  +1.40%  "

# Row 31 - Cosmos
$ws.Range("D31").Value = "9.07"
$ws.Range("E31").Value = "  -1.13%  "

# Row 32 - InjectiveProtocol
$ws.Range("D32").Value = "33.11"
$ws.Range("E32").Value = "  +0.58%  "

# Row 33 - FirstDigitalUSD
$ws.Range("E33").Value = "  +0.16%  "

# Row 34 - RenderToken
$ws.Range("D34").Value = "4.83"
$ws.Range("E34").Value = "  -0.30%  "

# Row 35 - Filecoin
$ws.Range("D35").Value = "5.02"
$ws.Range("E35").Value = "  -2.84%  "

# Row 36 - Celestia
$ws.Range("D36").Value = "17.27"
$ws.Range("E36").Value = "  -4.58%  "

# Row 37 - WEMIXToken (unchanged)

# Row 38 - Hedera
$ws.Range("D38").Value = "0.0693"
$ws.Range("E38").Value = "  +0.67%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -0.99%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -1.77%  "

# Row 41 - Stellar
$ws.Range("E41").Value = "  -1.16%  "

# Row 42 - LidoDAOToken
$ws.Range("E42").Value = "  -0.80%  "

# Row 43 - Maker
$ws.Range("D43").Value = "2.008.16"
$ws.Range("E43").Value = "  +0.11%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -2.12%  "

# Row 45 - FraxShare
$ws.Range("D45").Value = "10.02"
$ws.Range("E45").Value = "  -1.34%  "

# Row 46 - ApeXProtocol
$ws.Range("D46").Value = "2.11"
$ws.Range("E46").Value = "  -2.33%  "

# Row 47 - EnergySwap
$ws.Range("E47").Value = "  +3.76%  "

# Row 48 - NEARProtocol
$ws.Range("D48").Value = "2.79"
$ws.Range("E48").Value = "  -1.45%  "

# Row 49 / 50 - MultiversX and HuobiToken swap ranking positions
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "54.28"
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "2.84"
$ws.Range("E50").Value = "  +2.57%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.527.11"
$ws.Range("E51").Value = "  -0.01%  "
